$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 254.33333
$ws.Range("J33").Value = 200
$ws.Range("L33").Value = 200
$ws.Range("N33").Value = -658
$ws.Range("H58").Value = 1915.6154
$ws.Range("J58").Value = 3019.3333
$ws.Range("L58").Value = 9057.999899999999
$ws.Range("N58").Value = -9357.999899999999
$ws.Range("H74").Value = 4457.316
$ws.Range("J74").Value = 5270.5
$ws.Range("L74").Value = 5270.5
$ws.Range("N74").Value = -7142.5
$ws.Range("H77").Value = 4457.316
$ws.Range("J77").Value = 5270.5
$ws.Range("L77").Value = 26352.5
$ws.Range("N77").Value = -35712.5
$ws.Range("H112").Value = 66379.12
$ws.Range("J112").Value = 101872.27
$ws.Range("L112").Value = 305616.81
$ws.Range("N112").Value = -307832.81
$ws.Range("H129").Value = 71429944
$ws.Range("I129").Value = 912
$ws.Range("K129").Value = 2736
$ws.Range("M129").Value = 2264
$ws.Range("H137").Value = 315848.38
$ws.Range("I137").Value = 408341.2
$ws.Range("K137").Value = 1225023.6
$ws.Range("M137").Value = -1222473.6
$ws.Range("H138").Value = 4338.2183
$ws.Range("J138").Value = 5090.2114
$ws.Range("L138").Value = 15270.6342
$ws.Range("N138").Value = -25550.6342

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 519000000
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H63").Value = 6754.364
$ws.Range("I63").Value = 6754.364
$ws.Range("K63").Value = 6754.364
$ws.Range("M63").Value = -6068.364
$ws.Range("H66").Value = 6754.364
$ws.Range("I66").Value = 6754.364
$ws.Range("K66").Value = 33771.82
$ws.Range("M66").Value = -30339.82
$ws.Range("H80").Value = 70851.03
$ws.Range("J80").Value = 100000
$ws.Range("L80").Value = 100000
$ws.Range("N80").Value = -101996
$ws.Range("H83").Value = 70851.03
$ws.Range("J83").Value = 100000
$ws.Range("L83").Value = 300000
$ws.Range("N83").Value = -309984
$ws.Range("H95").Value = 81500
$ws.Range("J95").Value = 81500
$ws.Range("L95").Value = 81500
$ws.Range("N95").Value = -86992
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H99").Value = 519000000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 100377
$ws.Range("J109").Value = 100377
$ws.Range("L109").Value = 100377
$ws.Range("N109").Value = -103151
$ws.Range("H110").Value = 18718.5
$ws.Range("I110").Value = 22062.2
$ws.Range("K110").Value = 22062.2
$ws.Range("M110").Value = -20017.2
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 58550
$ws.Range("J106").Value = 58550
$ws.Range("L106").Value = 58550
$ws.Range("N106").Value = -61074
$ws.Range("H107").Value = 7673.75
$ws.Range("I107").Value = 7785.3335
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 7785.3335
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = -5865.3335
$ws.Range("N107").Value = -9840
$ws.Range("H108").Value = 455703330
$ws.Range("J108").Value = 455703330
$ws.Range("L108").Value = 455703330
$ws.Range("N108").Value = -455711010
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H117").Value = 742000
$ws.Range("J117").Value = 742000
$ws.Range("L117").Value = 742000
$ws.Range("N117").Value = -751178
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H140").Value = 89299
$ws.Range("J140").Value = 89299
$ws.Range("L140").Value = 89299
$ws.Range("N140").Value = -99659

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 645.44446
$ws.Range("I7").Value = 827.8
$ws.Range("J7").Value = 124.42857
$ws.Range("K7").Value = 827.8
$ws.Range("L7").Value = 124.42857
$ws.Range("M7").Value = -714.8
$ws.Range("N7").Value = -350.42857
$ws.Range("H22").Value = 722.0769
$ws.Range("I22").Value = 733.1667
$ws.Range("J22").Value = 712.5714
$ws.Range("K22").Value = 733.1667
$ws.Range("L22").Value = 712.5714
$ws.Range("M22").Value = -383.1667
$ws.Range("N22").Value = -1412.5714
$ws.Range("H131").Value = 777
$ws.Range("J131").Value = 777
$ws.Range("L131").Value = 777
$ws.Range("N131").Value = -10857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 404.375
$ws.Range("J92").Value = 383
$ws.Range("L92").Value = 1149
$ws.Range("N92").Value = -3645

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 19050
$ws.Range("I43").Value = 19050
$ws.Range("K43").Value = 19050
$ws.Range("M43").Value = -18899
$ws.Range("H57").Value = 44969
$ws.Range("J57").Value = 44969
$ws.Range("L57").Value = 44969
$ws.Range("N57").Value = -46609
$ws.Range("H80").Value = 4993.3335
$ws.Range("I80").Value = 4993.3335
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4993.3335
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3995.3335
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 4993.3335
$ws.Range("I83").Value = 4993.3335
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 24966.6675
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -19974.6675
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 10079.481
$ws.Range("I97").Value = 12653.4
$ws.Range("K97").Value = 12653.4
$ws.Range("M97").Value = -12157.4
$ws.Range("H122").Value = 12799.517
$ws.Range("I122").Value = 11517.88
$ws.Range("J122").Value = 18139.666
$ws.Range("K122").Value = 34553.64
$ws.Range("L122").Value = 54418.99800000001
$ws.Range("M122").Value = -32103.64
$ws.Range("N122").Value = -59318.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26573.666
$ws.Range("I7").Value = 31519.21
$ws.Range("K7").Value = 31519.21
$ws.Range("M7").Value = -31407.21
$ws.Range("H22").Value = 4221
$ws.Range("I22").Value = 4221
$ws.Range("K22").Value = 4221
$ws.Range("M22").Value = -3926
$ws.Range("H27").Value = 4221
$ws.Range("I27").Value = 4221
$ws.Range("K27").Value = 4221
$ws.Range("M27").Value = -4114
$ws.Range("H46").Value = 2770.1428
$ws.Range("J46").Value = 4664.3335
$ws.Range("L46").Value = 4664.3335
$ws.Range("N46").Value = -5040.3335
$ws.Range("H126").Value = 26573.666
$ws.Range("I126").Value = 31519.21
$ws.Range("K126").Value = 94557.63
$ws.Range("M126").Value = -92087.63

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 31679.7
$ws.Range("J41").Value = 43716.168
$ws.Range("L41").Value = 43716.168
$ws.Range("N41").Value = -44496.168
$ws.Range("H45").Value = 21804
$ws.Range("J45").Value = 21804
$ws.Range("L45").Value = 21804
$ws.Range("N45").Value = -22786
$ws.Range("H74").Value = 37631.6
$ws.Range("I74").Value = 5466.6665
$ws.Range("J74").Value = 45672.832
$ws.Range("K74").Value = 5466.6665
$ws.Range("L74").Value = 45672.832
$ws.Range("M74").Value = -4530.6665
$ws.Range("N74").Value = -47544.832
$ws.Range("H77").Value = 37631.6
$ws.Range("I77").Value = 5466.6665
$ws.Range("J77").Value = 45672.832
$ws.Range("K77").Value = 16399.9995
$ws.Range("L77").Value = 137018.496
$ws.Range("M77").Value = -11719.9995
$ws.Range("N77").Value = -146378.496
$ws.Range("H107").Value = 2545.9375
$ws.Range("I107").Value = 2284
$ws.Range("K107").Value = 6852
$ws.Range("M107").Value = -4932
$ws.Range("H126").Value = 32125.867
$ws.Range("I126").Value = 53612.25
$ws.Range("J126").Value = 7570
$ws.Range("K126").Value = 160836.75
$ws.Range("L126").Value = 22710
$ws.Range("M126").Value = -158366.75
$ws.Range("N126").Value = -27650
